$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G22").Value = 'Comentarios:
- Lograste desarrollar un gran proyecto a pesar de las dificultades que tuviste a lo largo del semestre por circunstancias diferentes a lo académico. Felicitaciones!
- Aquí van algunos comentarios que pueden servirte de cara al reto del siguiente semestre: 1) el tamaño de la cocina puede ser muy reducido para lo que se necesita en esta: es importante que haya separación entre la zona de cocción y de lavado; 2) Eventualmente, podrías sacar la nevera para tener zona de preparación; 3) en la medida de lo posible, es importante que trates de lograr que todos los baños cuenten con ventilación natural; 4) para que la zona de ropas funcione mejor (que la ropa se sequé de manera adecuada), es importante que esta cuente con buena ventilación.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,6. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,6.'
$ws.Range("G47").Value = 'Comentarios:
- Tuviste un gran proceso a lo largo del semestre. Además del aprendizaje netamente disciplinar de la arquitectura, nos gustó mucho ver cómo ganaste en confianza con relación a ti misma. Tienes muy buenas bases para el reto que encontrarás en el siguiente taller.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) trata de conciliar mejor la lógica constructiva con los espacios que diseñas; 2) en una cultura como la de estos habitantes, es importante contemplar la presencia de un comedor.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G72").Value = 'Comentarios:
- Valoramos el esfuerzo por terminar  dentro del plazo establecido este capítulo en tu proceso de formación, pese al reto personal que estás enfrentando: mucha fuerza. 
- Te invitamos a seguir desarrollando las capacidades con que cuentas, con mayor decisión y constancia a lo largo de todo el proceso y no solo en momentos puntuales de una entrega.
- En un próximo proyecto de vivienda, es fundamental que todas las habitaciones puedan ventilarse e iluminarse de manera natural (esto no sucede con la habitación de la persona con discapacidad).
- La nota de tu entrega final, correspondiente al 30% del curso es 3,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,8.'
$ws.Range("G97").Value = 'Comentarios:
- Tienes grandes capacidades y estás preparada para asumir retos de mayor complejidad. Logras un proyecto de buen nivel. 
- Esperamos que puedas gestionar mejor los momentos de entrega, evitando tener rebajas por no hacerlo a tiempo. 
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) para tener un lote tan grande, pudiste haber desarrollado un diseño más complejo y rico; 2) el baño compartido de las habitaciones se encuentra muy alejado de estas; 3) no hay distribución de los elementos de cocina (nevera, estufa y lavaplatos) y esto es importante para reconocer los flujos de las actividades.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,3, sin embargo, al hacerloen un día posterior, se computa en 2,3. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,6.'
$ws.Range("G122").Value = 'Comentarios:
- Estás preparado para desarrollar proyectos de mayor complejidad. 
- Te invitamos a no dejar de usar tus manos para modelar maquetas: esto te permite entender mejor el espacio y evaluar alternativas de manera más versátil.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) es fundamental representar los elementos de la cocina para evaluar qué tan adecuadamente se pueden llevar a cabo las actividades en su interior (nevera, estufa y lavaplatos); 2) la zona de ropas debería tener ventilación cruzada en un lugar como Honda, el cual tiene una alta humedad relativa; 3) baño con 2 accesos requiere resolver mejor la privacidad de quien lo usa.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4.'
$ws.Range("G147").Value = 'Comentarios:
- Estás preparado para desarrollar un proyecto de mayor complejidad. 
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) a pesar de contar con un lote tan generoso, la habitación 1 tiene unas dimesiones mínimas; 2) por el contrario, los baños son más grandes de lo convencional sin una intención que lo justifique; 3) la cocina tiene pocos espacios de almacenaje: para esto puede ser útil hacer un análisis de uso.
- La nota de tu entrega final, correspondiente al 30% del curso es 4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4.'
$ws.Range("G172").Value = 'Comentarios:
- Tienes grandes capacidades que pudimos evidenciar a lo largo de todo el proceso del taller y se demuestran en la entrega que presentas. No dejes de hacer búsquedas conectadas con tus intereses particulares, que van más allá de lo que se desarrolla a nivel general en un taller de diseño: esto te da un diferencial. Te invitamos a seguir trabajando el reto que tienes con el perfeccionismo, pues se puede convertir en tu mayor amenaza.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) conviene que la cocina sea un poco más amplia, especialmente para tener mayor área de preparación; 2) una distribución en L para este caso, podría ayudarte; 3) no aparecen representadas las zonas de almacenaje en habitaciones (closet´s).
- La nota de tu entrega final, correspondiente al 30% del curso es 4,9. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,7.'
$ws.Range("G197").Value = 'Comentarios:
- Valoramos la manera como relacionas el interior con el exterior de tu vivienda en el costado más largo del lote, tratando de ofrecer algo hacia el municipio de Honda y no solamente para tu familia. 
- Estás preparado para desarrollar proyectos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) no es claro el lugar en el cual la familia puede compartir el ritual de comer, en medio de diferentes condiciones climáticas... esto es importante para esta cultura.; 2) el corredor de las habitaciones no tiene iluminación y conviene explorar estrategias para lograrlo.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,1.'
$ws.Range("G222").Value = 'Comentarios:
- No pierdas el espíritu combativo para sobreponerte a las dificultades que encuentras. Es necesario que hagas esfuerzos adicionales para terminar de aprender lo que exige este taller: esto es fundamental para poder gestionar la complejidad mayor que encontrarás en el siguiente semestre. Esperamos que la fuerza que tuviste para tomar la decisión personal de continuar al siguiente nivel, la enfoques en ejercicios autogestionados durante el período intersemestral; estos podrían ir en tres frentes: 1)  entender mejor el espacio en sus 3 dimensiones y no en representaciones bidimensionales aisladas; 2) diseñar una estructura que haga factible las intenciones espaciales que tienes; 3) afianzar el dominio técnico de las herramientas de representación y modelado.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) en la cocina no están representados todos los elementos (nevera, estufa y lavaplatos), por lo cual es dificil evaluar qué tan bien funciona este espacio en función de las actividades que se llevan a cabo en su interior; 2) convendría tener un baño privado en la zona principal; 3) la propuesta estructural de la casa en el primer y segundo nivel no es clara y hacen falta elementos de soporte que la hagan factible; 4) no se entiende la manera en que se soportan las escaleras; 5) hay intenciones espaciales en planta que se desdibujan con lo que finalmente desarrollas en las fachadas.
- La nota de tu entrega final, correspondiente al 30% del curso es 3,6; al hacerlo por fuera del momento fijado, tienes una rebaja  de una unidad y se computa en 2,6. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,6.'
$ws.Range("G247").Value = 'Comentarios:
- Llevaste a cabo un proceso caracterizado por la constancia, la disciplina y la puntualidad. Presentas un proyecto bien desarrollado y estás lista para abordar retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) en la zona de lavado se puede pensar en un espacio más demarcado; 2) convien prever el progreso de la familia (¿qué pasaría si cambian de vehículo?)
- La nota de tu entrega final, correspondiente al 30% del curso es 4,9. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G272").Value = 'Comentarios:
- Tuviste un proceso impecable y estás lista para asumir nuevos retos.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) no es lógico que el acceso para la habitación principal deba hacerse a través de la zona de lavado; 2) el baño social está sobredimensiado.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G297").Value = 'Comentarios:
- Muy valiosos los argumentos que te llevaron a tomar decisiones en tu vivienda (presentaste una gran narrativa). Los conocimientos adquiridos te habilitan para desarrollar retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) el baño de la habitación principal no tiene acceso de acuerdo a lo que se ve en la planta; 2) el tamaño de éste es muy reducido, el cual aunque cumple con estándares mínimos, podría ser más generoso teniendo presente el tamaño del lote; 3) conviene reflexionar sobre las actividades y la ubicación de acceso a algunos espacios (ej: habitación principal)
- La nota de tu entrega final, correspondiente al 30% del curso es 4,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,7.'
$ws.Range("G322").Value = 'Comentarios:
- Estás preparada para desarrollar retos de mayor complejidad. Que las restricciones que encuentres al pensar en la vivienda colectiva, potencien tu creatividad y habilidades. 
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) aunque el área de la cocina es generosa, es importante relacionar la distribución de elementos (estufa, lavaplatos, etc) para visualizar las posibilidades de moverse en el espacio; 2) las zonas húmedas (ej: piscina) deben contemplar baños cerca; 3) el baño social está sobredimensionado; 4) al tener un lote tan grande, con una fachada tan extensa hacia el espacio público de Honda, es conveniente ofrecer algo más rico que un muro (pensar la casa no solo de puertas para adentro, sino en qué le aporta esta al entorno cercano).
- La nota de tu entrega final, correspondiente al 30% del curso es 3,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,4.'
$ws.Range("G347").Value = 'Comentarios:
- Es muy valioso cómo, a partir de tu observación y sensibilidad, lograste relacionar las dinámicas de un gran río con la experiencia vivida en la salida de campo, para tomar decisiones en tu proyecto que consideran espacios para personas en tránsito. Estás preparada para asumir retos más complejos.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) el baño del segundo piso pudo ser más amplio; 2) convendría explorar alternativas que le permitan a la familia estar más cerca... en nuestro contexto, los hijos suelen ubicarse cerca de la habitación principal; 3) no hay diferencias claras que demuestren cuál es la habitación principal.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,6. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,7.'
$ws.Range("G372").Value = 'Comentarios:
- Interesante el desafío particular que te formulaste, desarrollando una vivienda que tiene una actividad productiva. Estás preparado para desarrollar retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) el baño social no debería tener el acceso justo al frente del comedor; 2) la cocina no tiene zona para preparar alimentos; 3) convendría explorar la posición del baño al fondo y tener la sala integrada con el comedor.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,8.'
$ws.Range("G397").Value = 'Comentarios:
- Decidiste emplear una geometría no convencional para articular todo el proyecto, lo que implicaba un reto de gran complejidad que desarrollaste con muy buen nivel. Lograste conectar intenciones plásticas con un desarrollo técnico que las hace viables, por ejemplo, la unión entre muros y cubierta, así como la claridad para definir las superficies en que se fragmenta una cubierta en espiral. Tienes una gran energía y estás listo para los retos que vienen.
- Para un próximo proyecto de vivienda: los baños requieren iluminación y ventilación (no se evidencia en la representación planimétrica); 2) el tamaño de las habitaciones versus la posición y posibilidades de los muebles, pueden tener ajustes.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,9. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G422").Value = 'Comentarios:
- Valoramos tu disposición para explorar alternativas en un lote con muchas restricciones. Fuiste recursivo y apasionado para desarrollar los diferentes laboratorios y esto te permite alcanzar un resultado muy bueno. Esperamos que el buen dominio que estás adquiriendo de herramientas digitales, no anulen la gran capacidad que tienes con dibujos hechos a mano.
- En un próximo proyecto de vivienda, la zona de ropas podría tener calados hacia la calle para que funcione mejor y el baño social podría ventilarse mejor.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,6.'
$ws.Range("G447").Value = 'Comentarios:
- Desarrollas un buen proyecto y estás preparado para elevar el complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) es importante representar todos los elementos de una cocina para analizar los movimientos humanos  (nevera, estufa y lavaplatos); 2) en un clima como el de Honda, con una humedad relativa tan alta, es fundamental separar los closet´s del baño.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,7. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,6.'
$ws.Range("G472").Value = 'Comentarios:
- Llevas el proyecto a un buen nivel de desarrollo y estás preparado parauna mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) el ingreso al baño del taller debió hacerse al interior de este; 2) hay una desproporción entre el tamaño de los baños y las habitaciones; 3) el comentario anterior también aplica para el vestidor.
- La nota de tu entrega final, correspondiente al 30% del curso es 4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,8.'
$ws.Range("G497").Value = 'Comentarios:
- Alcanzas un buen resultado y estás preparada para desarrollar proyectos más complejos.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) en la cocina debe separarse la zona húmeda de la zona de cocción; 2) no es clara la zona de almacenaja en la habitación principal (closet); 3) una cubierta translúcida sobre un espacio interno, en un clima como este no es una buena opción.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,5. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,8.'
$ws.Range("G522").Value = 'Comentarios:
- Desarrollas el proyecto a un buen nivel y estás preparada para incremetar el nivel de complejidad. Te invitamos a confiar más en tí y no pretender que la opinión de todas las personas a las que les haces preguntas coincidan.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) la cocina del restaurante debería estar mejor ventilada; 2) conviene tener un cerramiento en la zona de ropas, de tal manera que cuando se esté usando (haya ropa extendida), no la vea todo el mundo; 3) no siempre contarás con áreas tan generosas, por lo tanto es importante reflexionar sobre la ubicación del baño social; 4) es fundamental para lo que viene de aquí en adelante, desarrollar estructuras en sus columnas (al menos la mayoría) tengan continuidad hasta el suelo.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,3. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,4.'
$ws.Range("G547").Value = 'Comentarios:
- Aprovecha toda la energía y fortaleza con la que cuentas, para hacer los esfuerzos adicionales que requieres para resolver los retos del próximo taller: además de ejercitarte en el dominio de herramientas de representación, es importante que seas capaz de desarrollar intenciones conceptuales al nivel técnico que requieren para que sean posibles. 
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) reflexionar sobre la ubicación de los accesos, especialmente el del baño que queda frente a la zona social; 2) la distribución de la cocina debe contemplar tres áreas (cocción, lavado y refrigeración); 3) el almacenaje es fundamental.
- La nota de tu entrega final, correspondiente al 30% del curso es 3,3. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,9.'
$ws.Range("G572").Value = 'Comentarios:
- Valoramos la persistencia y la disciplina que tuviste a lo largo del proceso, resolviendo una vivienda con una actividad productiva con requerimientos de un nivel superior a los que plantea el taller. Así mismo, el hecho de no haber renunciado a expresarte de manera oral en diferentes momentos del taller, hablan de tu fortaleza y tenacidad. Estás preparada para resolver proyectos más complejos.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) el baño del segundo piso queda frente a la zona social y esto no es conveniente (existen posibilidades de modificar el acceso); 2) al tener patios, la zona de lavado no tiene tanta lógica en el segundo piso.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,4.'
$ws.Range("G597").Value = 'Comentarios:
- Tienes una gran capacidad para vincular las experiencias vividas durante la salida a Honda con tus convicciones, tus propósitos y lo que estás aprendiendo en la profesión. Es muy valioso que hayas decidido hacer ajustes cuando el proyecto ya estaba avanzado, lo que te permitió desarrollarlo a un nivel muy bueno: ejercitaste una habilidad que será esencial en tu práctica profesional (hacer cambios a lo largo del proceso en función de lo que va emergiendo en el proceso). Finalmente, hiciste una exposición excelente de tu proyecto, tanto por la narrativa que preparaste, como por la manera en que conectaste con la audiencia que tenías.
- En un próximo proyecto de vivienda, la distribución de habitaciones podría integrarse más con los baños y las zonas de almacenaje de ropa. Como está en este momento, se ve como algo añadido que no pertenece al espacio.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,9. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G622").Value = 'Comentarios:
- Que este alto en tu proceso de formación y la energía que se libera al hacer la elección que tomaste, te impulsen a desarrollar ejercicios autónomos que te permitan crecer en tres aspectos: 1) tener un mayor dominio de herramientas de representación; 2) entender la lógica que hay detrás de aquello que representas, haciendo más habitable y factible lo que diseñas, además de ser consciente de los elementos que requieres a nivel constructivo para que esto sea posible; 3) usar el tiempo de manera efectiva durante las sesiones de taller y no perder el contacto con los profesores que están a cargo de este.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) la zona de ropas no tiene ventilación y esto lo convierte en un lugar húmedo con problemas para su funcionamiento; 2) la proporción de las habitaciones no corresponde a las actividades (hay una en particular que es muy angosta).
- La nota de tu entrega final, correspondiente al 30% del curso es 2,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,5.'
$ws.Range("G647").Value = 'Comentarios:
- Tienes grandes capacidades demostradas a lo largo de todo el proceso. Siempre te vimos dispuesta, atenta, concentrada y proactiva. Tuviste el valor de hacer ajustes que demandaban una energía extra, pese a tener un proyecto muy avanzado varias semanas antes de finalizar el proyecto: esto es algo muy importante para tu ejercicio profesional, en el que estarás expuesta a situaciones de este tipo.
- En un próximo proyecto de vivienda, la zona de ropas y el baño del primer piso se encuentran sobredimensionados. Por el contrario, la cocina quedó muy justa en sus dimensiones; eventualmente, una distribución de la cocina en "L" te permitiría ganar espacio para la zona de preparación.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G672").Value = 'Comentarios:
- Alcanzas a desarrollar un proyecto con un buen nivel y estás preparada para asumir retos de mayor nivel de complejidad. Pudimos notar el cambio en tu manera de trabajar, al ser consciente de la atención y dedicación que requerías para poder tener un cierre sin contratiempos.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) el acceso principal puede reducirse; 2) la solución de puertas corredizas debe revisarse, porque como están desarrolladas en este momento, hay conflictos con los muros portantes.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,5. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,2.'
$ws.Range("G697").Value = 'Comentarios:
- Alcanzas un resultado impecable: detalles como el simple hecho de variar el enrase de los muros con la cubierta en la maqueta, demuestran el rigor y la entrega que invertiste en tu proceso de aprendizaje (en la mayoría de proyectos esto no se hizo). Tienes una gran disciplina y es evidente el cariño que tienes por aquello que haces. Te invitamos a seguir exponiéndote con mayor apertura y decisión ante otras personas que no están tan cerca de ti.
- En un próximo proyecto de vivienda, conviene pensar mejor los accesos a los espacios: por ejemplo, acceder a la zona de ropas a través del comedor no es tan bueno, cuando esto podía hacerse a través del patio.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,9. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G722").Value = 'Comentarios:
- Desarrollas un buen proyecto y estás lista para resolver retos de mayor complejidad.
- En un próximo proyecto de vivienda, es bueno reducir un poco las áreas de circulación (están sobredimensionadas).
- La nota de tu entrega final, correspondiente al 30% del curso es 4,5. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,7.'
$ws.Range("G747").Value = 'Comentarios:
- Tuviste un gran proceso caracterizado por la atención, disciplina, fortaleza y la exploración. Muy valioso que hayas aceptado hacer cambios que implicaban reprocesos en un proyecto que estaba en una etapa avanzada. Así mismo, el hecho de tomar el riesgo por desarrollar una cubierta no convencional, que propone cosas interesantes para ese contexto.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) contemplar zonas de almacenaje dentro del mismo diseño (closet´s conformados por nichos entre muros, en lugar de muebles sobre una pared); 2) desarrollar con mayor detalle, la manera en que se encuentran los muros con la cubierta (es una cercha? son muros con enrases inclinados?)
- La nota de tu entrega final, correspondiente al 30% del curso es 4,7. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G772").Value = 'Comentarios:
- Presentas un gran proyecto fruto de un gran proceso caracterizado por el amor, la sensibilidad, dedicación y entrega. Cuentas con muy buenas bases para seguir adelante en tu proceso de aprendizaje.
- En un próximo proyecto de vivienda, los accesos de los baños no deberían quedar al frente de una zona social (ej: al frente del comedor).
- La nota de tu entrega final, correspondiente al 30% del curso es 5. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G797").Value = 'Comentarios:
- Esperamos que aprendas mucho en el proceso autónomo que inicias, de tal manera que puedas gestionar de manera adecuada la complejidad mayor que encontrarás en el siguiente taller. Puedes plantearte ejercicios concretos que respondan a 4 preguntas guía: ¿qué requieren ciertos espacios según su uso para hacerlos más habitables?; ¿cómo pueden relacionarse diferentes espacios entre sí de una manera más adecuada, respondiendo a las secuencias de uso de una vivienda?; ¿cómo mejorar el dominio técnico de las herramientas de representación?; ¿qué elementos constructivos se requieren para hacer posible lo que estás diseñando?. Finalmente, te invitamos a usar el tiempo de manera más efectiva durante las sesiones de taller y no perder contacto en momentos cruciales con los profesores del taller.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) algunos espacios se encuentran desproporcionados para el uso que tienen (ej: baño con 2 duchas (cosa que no entendimos) es casi del mismo tamaño que la habitación auxiliar); 2) corredor tiene un ancho excesivo; 3) la tienda por norma debe tener baño debido al hecho de tener alimentos; 4) los accesos a los baños y a los closet´s son muy pequeños ´a pesar de ser muy generosos en área, en algunos casos tienen problemas para acceder por la ubicación de elementos como un lavamanos doble; 5) a pesar de la gran área del lote, hay baños sin iluminación ni ventilación natural; 6) corredor remate en la puerta de un baño.
- La nota de tu entrega final, correspondiente al 30% del curso es 2,5. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3.'
$ws.Range("G822").Value = 'Comentarios:
- Desarrollas tu proyecto a un buen nivel y estás lista para retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) hay circulaciones que hacen perder áreas (zona de lavado por ejemplo); 2) la cocina es generosa, pero no está clara la distribución de las zonas que la componen (cocción, refrigeración, lavado y preparación); 3) la habitación 4 no tiene un baño cerca, implicando un recorrido excesivo para los habitantes; 4) el acceso de la vivienda no tiene las medidas mínimas.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,2. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,3.'
$ws.Range("G847").Value = 'Comentarios:
- Tuviste un proceso equilibrado a lo largo del taller y esto te permitió desarrollar un proyecto interesante y completo. Estás lista para un reto de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) convendría explorar qué sucede si la cocina se diseña en "L", con el fin de ganar área de preparación y disminuir áreas de circulación; 2) el acceso al baño social quedó frente al comedor (puede reubicarse); 3) zona de ropas muy encerrada, pese a tener grandes posibilidades para ventilarse.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,2. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,3.'
$ws.Range("G872").Value = 'Comentarios:
- El resultado que alcanzas demuestra que estás lista para abordar retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) a pesar de tener un lote generoso, puede estarte haciendo falta un baño en el segundo piso, porque el baño destinado a los invitados, en realidad es el de las habitaciones; 2) es importante presentar material actualizado según la propuesta: el patio del boceto no corresponde al que finalmente tiene la vivienda (revisar proporciones).
- La nota de tu entrega final, correspondiente al 30% del curso es 4,5. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,7.'
$ws.Range("G897").Value = 'Comentarios:
- Fuiste muy recursiva para resolver una vivienda en un lote con muchas restricciones. Desde el inicio hasta el fin, tu proceso se caracterizó por una gran dedicación, atención y disposición. Estás lista para abordar retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) una de las habitaciones es demasiado pequeña; 2) el baño social y la zona de ropas podrían estar mejor ventilados a través de la cubierta por ejemplo.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,5.'
$ws.Range("G922").Value = 'Comentarios:
- Las capacidades que tienes, además de la actitud que tuviste a lo largo de un proceso constante, te permiten alcanzar un resultado muy bueno. Estás lista para resolver retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) en la distribución de la cocina, no se representan todos los elementos  (nevera, estufa y lavaplatos); 2) se podrían aprovechar mejor los ángulos, para ampliar las superficies de preparación; 3) el baño social debería tener su acceso en otro lugar: quedó directo a la zona social.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,7. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,6.'
$ws.Range("G947").Value = 'Comentarios:
- Has llevado a cabo un proceso muy bueno a lo largo del taller y alcanzas un buen desarrollo en lo que presentas en la entrega final. Vemos que estás creciendo en madurez y confianza contigo mismo: por ahí es.
- En un próximo proyecto de vivienda, conviene que se integren mejor los elementos de almacenaje en las habitaciones (closet´s).
- La nota de tu entrega final, correspondiente al 30% del curso es 4,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,7.'
$ws.Range("G972").Value = 'Comentarios:
- Valoramos el hecho de haberte retado con una geometría no convencional que exigía un esfuerzo extra. Estás lista no solo para resolver retos de diseño más complejos, sino también para enfrentar situaciones propias del ejercicio profesional gracias a tu gran habilidad para relacionarte con otros.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) convendría que un módulo contuviera el baño social y la zona de ropas, de tal manera que su acceso no quede expuesto a la sala; 2) la humedad puede traer problemas para los closet´s al estar ubicados al interior del baño.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,2. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,1.'
$ws.Range("G997").Value = 'Comentarios:
- Presentas un proyecto con un nivel muy bueno y estás preparada para desarrollar retos de mayor complejidad. Fue muy valioso ver la disposición y apertura que mostraste a lo largo del proceso para responder a lo que ibas descubriendo. Esto se evidenció especialmente cuando realizaste cambios que permitieron adaptar el proyecto a la topografía del lote. Las grandes habilidades que tienes a nivel personal y la actitud con que asumes tu vida te llevarán muy adelante.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) los patios pueden aprovecharse mejor para iluminar y ventilar los espacios contiguos a estos; 2) es posible ampliar un poco la zona de ropas y abrirla más hacia el talud; 3) el baño del segundo piso requiere más iluminación.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,8.'
$ws.Range("G1022").Value = 'Comentarios:
- Presentas un proyecto muy interesante que propone alternativas fuera de lo convencional y que desarrollas a un nivel muy bueno.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) los muebles no están bien escalados (revisar baños y cocina); 2) en la cocina es importante separar la zona húmeda de la de cocción; 3) los baños requieren ventilación.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,8. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,9.'
$ws.Range("G1047").Value = 'Comentarios:
- Tienes grandes habilidades que te permiten continuar en tu proceso de aprendizaje, abordando retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) considerar cuál es el tamaño mínimo para una habitación (hay una muy pequeña); 2) integrar mejor los baños; 3) en el segundo piso puede estar sobrando un baño; 4) si la cocina estuviera en L, podría aprovecharse mejor el espacio disponible (tener más área para preparar alimentos).
- La nota de tu entrega final, correspondiente al 30% del curso es 4,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,5.'
$ws.Range("G1072").Value = 'Comentarios:
- Valoramos la decisión que tomaste en las sesiones finales del taller, con el fin de hacer ajustes que mejoraban la experiencia que tendrían futuros habitantes de ese espacio. Estás preparado para desarrollar retos de mayor complejidad.
- En un próximo proyecto de vivienda, es extraño tener un corredor en el segundo piso que lleva a un "mini balcon". En su lugar, este espacio podría haberse destinado para closet´s entre las habitaciones, dándoles más espacio, más independencia y mayor tamaño para los balcones de cada una. 
- La nota de tu entrega final, correspondiente al 30% del curso es 3,9. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,2.'
$ws.Range("G1097").Value = 'Comentarios:
- Cumples con lo que se requiere en este nivel de formación y estás preparada para resolver retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) es fundamental representar los elementos de una cocina (nevera, estufa y lavaplatos), pues sin estos es difícil determinar la manera en que se realizarán las actividades; 2) considerar alternativas para ventilar los baños.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,2. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 3,8.'
$ws.Range("G1122").Value = 'Comentarios:
- El rigor, la constancia y la lógica fueron tus mejores aliados. Estás preparado para resolver retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) en las habitaciones no es clara la zona de almacenaje (closet´s); 2) podría tener más lógica ubicar el patio de ropas en el primer nivel: en esta cultura, la ropa se suele secar en el aire; 3) incorporar una fachada completamente vidriada, debería corresponder con una intención contundente hacia el paisaje; sin embargo, en este caso lo que se está viendo es un muro medianero. Esto aumenta costos de manera innecesaria, pero además, incrementa la ganancia de calor por cuenta del sol que ingresa al espacio, cosa que no conviene para nada en este clima.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,4. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,3.
'
$ws.Range("G1147").Value = 'Comentarios:
- Alcanzas a desarollar tu proyecto a un nivel muy bueno y estás preparada para abordar retos de mayor complejidad.
- Aquí van algunos comentarios que pueden servirte en un próximo proyecto de vivienda: 1) resolver mejor la ubicación de los baños; en el caso del baño del segundo piso, su acceso se encuentra frente a la zona social; 2) el baño de la tienda tiene la misma dificultad: está al frente de las mesas. Es importante pensar en la intimidad de las personas.
- La nota de tu entrega final, correspondiente al 30% del curso es 4,6. Teniendo presente las cuatro entregas previas con sus respectivos porcentajes, la nota general del curso que se subirá al sistema de la Universidad es 4,4.'
